$d = $word.ActiveDocument

# --- Change 1: merge "(Supreme Court)" heading runs into one run, drop bookmark ---
$d.Content.Find.Execute("(Supreme Court)", $true, $false, $false, $false, $false, $true, 1, $false, "(Supreme Court)", 2)

# --- Change 2: "supreme clerk." -> "Supreme Court clerk." (first occurrence, non-bold) ---
$d.Content.Find.Execute("supreme clerk.", $true, $false, $false, $false, $false, $true, 1, $false, "Supreme Court clerk.", 2)

# --- Change 3: "supreme clerk." -> "Supreme Court clerk." (second occurrence, bold) ---
$d.Content.Find.Execute("supreme clerk.", $true, $false, $false, $false, $false, $true, 1, $false, "Supreme Court clerk.", 2)
